$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 (bold, centered, bordered) onto the
# new H1 "Save" header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"
$excel.CutCopyMode = 0

# Populate column H ("Save") for each data row based on the "sum" column
# (G): rows with a sum of 9 or more are flagged as a save (1), others 0.
for ($r = 2; $r -le 48; $r++) {
    $sum = [double]$ws.Cells.Item($r, 7).Value2
    if ($sum -ge 9) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
